$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E16").Value = "2502"
$ws.Range("E17").Value = "2503"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2505"
$ws.Range("E20").Value = "2506"
$ws.Range("E21").Value = "2507"

$ws.Range("G16").Value = 1644396
$ws.Range("G17").Value = 1644396
$ws.Range("G18").Value = 1644396
$ws.Range("G19").Value = 1644396
$ws.Range("G20").Value = 1644396
$ws.Range("G21").Value = 1644396
